# Remove the blank "Sheet1" worksheet, then rename the remaining
# "test01" sheet to "Sheet1" and move the selection to A3.

$wb = $excel.ActiveWorkbook

# Suppress the "delete sheet" confirmation prompt.
$excel.DisplayAlerts = $false

$blankSheet = $wb.Worksheets.Item("Sheet1")
$blankSheet.Delete()

$excel.DisplayAlerts = $true

$ws = $wb.Worksheets.Item("test01")
$ws.Name = "Sheet1"

$ws.Activate()
$ws.Range("A3").Select()
